$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 10:52"


# Row 17
$ws.Range("B17").Value = 12058
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3463
$ws.Range("E17").Value = 8375
$ws.Range("F17").Value = 250
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 220

# Row 24
$ws.Range("B24").Value = 5795
$ws.Range("C24").Value = 45
$ws.Range("E24").Value = 3441

# Row 32
$ws.Range("E32").Value = 3333
$ws.Range("G32").Value = 6
$ws.Range("H32").Value = 157

# Row 34
$ws.Range("A34").Value = "Filipinas"
$ws.Range("B34").Value = 3660
$ws.Range("C34").Value = 414
$ws.Range("D34").Value = 73
$ws.Range("E34").Value = 3424
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 11
$ws.Range("H34").Value = 163

# Row 35
$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 3654
$ws.Range("D35").Value = 575
$ws.Range("E35").Value = 2994
$ws.Range("F35").Value = 69
$ws.Range("H35").Value = 85

# Row 36
$ws.Range("A36").Value = "Ecuador"
$ws.Range("B36").Value = 3646
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 100
$ws.Range("E36").Value = 3366
$ws.Range("F36").Value = 100
$ws.Range("H36").Value = 180

# Row 39
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 2491
$ws.Range("C39").Value = 218
$ws.Range("D39").Value = 192
$ws.Range("E39").Value = 2090
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 209

# Row 40
$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("B40").Value = 2463
$ws.Range("C40").Value = 61
$ws.Range("D40").Value = 488
$ws.Range("E40").Value = 1941
$ws.Range("F40").Value = 41
$ws.Range("H40").Value = 34

# Row 41
$ws.Range("A41").Value = "Peru"
$ws.Range("B41").Value = 2281
$ws.Range("D41").Value = 989
$ws.Range("E41").Value = 1209
$ws.Range("F41").Value = 81
$ws.Range("H41").Value = 83

# Row 66
$ws.Range("D66").Value = 37
$ws.Range("E66").Value = 810
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 17

# Row 71
$ws.Range("A71").Value = "Bielorrusia"
$ws.Range("C71").Value = 138
$ws.Range("D71").Value = 53
$ws.Range("E71").Value = 634
$ws.Range("F71").Value = 11
$ws.Range("G71").Value = 5
$ws.Range("H71").Value = 13

# Row 72
$ws.Range("A72").Value = "Barein"
$ws.Range("B72").Value = 700
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 431
$ws.Range("E72").Value = 265
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 4

# Row 73
$ws.Range("A73").Value = "Bosnia y Herzegovina"
$ws.Range("B73").Value = 667
$ws.Range("C73").Value = 13
$ws.Range("D73").Value = 44
$ws.Range("E73").Value = 597
$ws.Range("F73").Value = 4
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 26

# Row 74
$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 650
$ws.Range("D74").Value = 17
$ws.Range("E74").Value = 624
$ws.Range("F74").Value = 0
$ws.Range("H74").Value = 9

# Row 75
$ws.Range("A75").Value = "Azerbaiyan"
$ws.Range("D75").Value = 32
$ws.Range("E75").Value = 545
$ws.Range("F75").Value = 17
$ws.Range("H75").Value = 7

# Row 76
$ws.Range("A76").Value = "Kazajistan"
$ws.Range("B76").Value = 584
$ws.Range("D76").Value = 42
$ws.Range("E76").Value = 536
$ws.Range("F76").Value = 6
$ws.Range("H76").Value = 6

# Row 77
$ws.Range("A77").Value = "Tunez"
$ws.Range("B77").Value = 574
$ws.Range("D77").Value = 5
$ws.Range("E77").Value = 547
$ws.Range("F77").Value = 39
$ws.Range("H77").Value = 22

# Row 100
$ws.Range("A100").Value = "Estado de Palestina"
$ws.Range("B100").Value = 246
$ws.Range("C100").Value = 9
$ws.Range("D100").Value = 25
$ws.Range("E100").Value = 220
$ws.Range("F100").Value = 0
$ws.Range("H100").Value = 1

# Row 101
$ws.Range("A101").Value = "Vietnam"
$ws.Range("B101").Value = 241
$ws.Range("D101").Value = 95
$ws.Range("E101").Value = 146
$ws.Range("F101").Value = 8
$ws.Range("H101").Value = 0
